# Generate Report for Handoff
# Adds two new localization-status rows (for 62006172-...md and
# d4e04865-...md) to the "Overview" summary sheet and to the "de-de"
# per-language sheet. The "zh-cn" sheet is left untouched, since neither
# of the two new files has been localized into zh-cn yet.

$wb = $excel.ActiveWorkbook

$ovr  = $wb.Worksheets.Item("Overview")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1) "Overview" sheet - append two new rows (4 and 5) at the bottom of
#    the existing table.
# ---------------------------------------------------------------------

$ovrTable = $ovr.ListObjects.Item(1)

$row4 = $ovrTable.ListRows.Add()
$ovr.Range("A4").Value = "62006172-e8d9-42a7-84a1-8b6afec2c830.md"
$ovr.Range("B4").Value = "e2e\62006172-e8d9-42a7-84a1-8b6afec2c830.md"
$ovr.Range("C4").Value = ".md"
$ovr.Range("E4").Value = "N\A"
$ovr.Range("F4").Value = "Ready for handoff"
$ovr.Range("G4").Value = "2016-08-18 22:40:22"
$ovr.Range("B4").Style = "HyperLink"
$ovr.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$row5 = $ovrTable.ListRows.Add()
$ovr.Range("A5").Value = "d4e04865-4fad-421d-aeba-5f52de107caf.md"
$ovr.Range("B5").Value = "e2e\d4e04865-4fad-421d-aeba-5f52de107caf.md"
$ovr.Range("C5").Value = ".md"
$ovr.Range("E5").Value = "N\A"
$ovr.Range("F5").Value = "Ready for handoff"
$ovr.Range("G5").Value = "2016-08-18 22:40:22"
$ovr.Range("B5").Style = "HyperLink"
$ovr.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ovr.Hyperlinks.Add($ovr.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/62006172-e8d9-42a7-84a1-8b6afec2c830.md", "", "", "e2e\62006172-e8d9-42a7-84a1-8b6afec2c830.md") | Out-Null
$ovr.Hyperlinks.Add($ovr.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d4e04865-4fad-421d-aeba-5f52de107caf.md", "", "", "e2e\d4e04865-4fad-421d-aeba-5f52de107caf.md") | Out-Null

# ---------------------------------------------------------------------
# 2) "de-de" sheet - insert a new row for 62006172-...md right before
#    the existing 9c8f9aaa-...md row (keeping alphabetical/original
#    ordering), then append a new row for d4e04865-...md at the end.
# ---------------------------------------------------------------------

$dedeTable = $dede.ListObjects.Item(1)

# Remember the stale hyperlinks that live on the row we are about to
# push down, so we can re-create them at their new location afterwards.
$oldA3Address = $null
foreach ($h in $dede.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$3') { $oldA3Address = $h.Address }
}

# Push the existing row 3 (9c8f9aaa-...) down to row 4, freeing up row 3
# for the new 62006172-... entry, and carry the table formatting with it.
$dede.Rows("3").Insert()
$dedeTable.Resize($dede.Range("A1:P4"))

# Drop the stale hyperlink objects that did not follow the row shift
# (this engine keeps hyperlink anchors fixed on row-insert).
$toDelete = New-Object System.Collections.ArrayList
foreach ($h in $dede.Hyperlinks) {
    if (($h.Range.Address() -eq '$A$3') -or ($h.Range.Address() -eq '$I$3')) {
        $toDelete.Add($h) | Out-Null
    }
}
foreach ($h in $toDelete) { $h.Delete() }

# New row 3: 62006172-e8d9-42a7-84a1-8b6afec2c830.md
$dede.Range("A3").Value = "62006172-e8d9-42a7-84a1-8b6afec2c830.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "e2e"
$dede.Range("E3").Value = "ht"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "62006172-e8d9-42a7-84a1-8b6afec2c830.fc26c2ddb0379bb88fd9b684f4c42ab4b10aa805.de-de.xlf"
$dede.Range("H3").Value = "2016-08-18 22:40:22"
$dede.Range("K3").Value = "0001-01-01 00:00:00"
$dede.Range("M3").Value = "'True"
$dede.Range("O3").Value = "'False"
$dede.Range("A3").Style = "HyperLink"
$dede.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Re-create the hyperlink for the row that moved down to row 4
# (9c8f9aaa-...), pointing at its original target.
$dede.Hyperlinks.Add($dede.Range("A4"), $oldA3Address, "", "", "9c8f9aaa-3e95-482c-b984-0f564d7c39ea.md") | Out-Null

# New hyperlink for the freshly-inserted row 3.
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/62006172-e8d9-42a7-84a1-8b6afec2c830.md", "", "", "62006172-e8d9-42a7-84a1-8b6afec2c830.md") | Out-Null

# Append new row 5: d4e04865-4fad-421d-aeba-5f52de107caf.md
$newRow = $dedeTable.ListRows.Add()
$dede.Range("A5").Value = "d4e04865-4fad-421d-aeba-5f52de107caf.md"
$dede.Range("B5").Value = ".md"
$dede.Range("C5").Value = "Ready for handoff"
$dede.Range("D5").Value = "e2e"
$dede.Range("E5").Value = "ht"
$dede.Range("F5").Value = "'False"
$dede.Range("G5").Value = "d4e04865-4fad-421d-aeba-5f52de107caf.f038b7b73ce57e924c06c432128d9752db7bb291.de-de.xlf"
$dede.Range("H5").Value = "2016-08-18 22:40:22"
$dede.Range("K5").Value = "0001-01-01 00:00:00"
$dede.Range("M5").Value = "'True"
$dede.Range("O5").Value = "'False"
$dede.Range("A5").Style = "HyperLink"
$dede.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$dede.Hyperlinks.Add($dede.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d4e04865-4fad-421d-aeba-5f52de107caf.md", "", "", "d4e04865-4fad-421d-aeba-5f52de107caf.md") | Out-Null
